$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($cellRef, $val)
    $c = $ws.Range($cellRef)
    $origStyle = $c.Style
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = $origStyle
}

Set-TextValue "D2" "42.960.03"
Set-TextValue "D3" "2.305.20"
Set-TextValue "E3" "  +0.17%  "
Set-TextValue "E4" "  -0.14%  "
Set-TextValue "D5" "303.92"
Set-TextValue "E5" "  +1.18%  "
Set-TextValue "D6" "97.10"
Set-TextValue "E6" "  -0.23%  "
Set-TextValue "D7" "0.503"
Set-TextValue "E7" "  -2.01%  "
Set-TextValue "E8" "  -0.05%  "
Set-TextValue "D9" "0.501"
Set-TextValue "E9" "  -0.98%  "
Set-TextValue "D10" "35.43"
Set-TextValue "E10" "  -1.01%  "
Set-TextValue "E11" "  -0.05%  "
Set-TextValue "D12" "18.78"
Set-TextValue "E12" "  +4.96%  "
Set-TextValue "E13" "  +1.66%  "
Set-TextValue "D14" "6.90"
Set-TextValue "E14" "  +1.92%  "
Set-TextValue "D15" "2.664.98"
Set-TextValue "E15" "  +0.35%  "
Set-TextValue "D16" "2.320.94"
Set-TextValue "E16" "  +0.64%  "
Set-TextValue "E17" "  +0.71%  "
Set-TextValue "D18" "42.844.91"
Set-TextValue "E18" "  -0.08%  "
Set-TextValue "E19" "  -1.39%  "
Set-TextValue "D20" "0.0₃0898"
Set-TextValue "E20" "  -0.78%  "
Set-TextValue "E21" "  -0.17%  "
Set-TextValue "D22" "67.77"
Set-TextValue "E22" "  -0.16%  "
Set-TextValue "D23" "237.06"
Set-TextValue "E23" "  -1.56%  "
Set-TextValue "E24" "  +1.24%  "
Set-TextValue "E25" "  +0.12%  "
Set-TextValue "E26" "  -0.20%  "
Set-TextValue "D27" "24.90"
Set-TextValue "E27" "  -2.19%  "
Set-TextValue "D28" "2.38"
Set-TextValue "E28" "  +17.77%  "
Set-TextValue "D29" "165.94"
Set-TextValue "E29" "  +0.13%  "
Set-TextValue "D30" "9.06"
Set-TextValue "E30" "  +0.21%  "
Set-TextValue "D31" "32.83"
Set-TextValue "E31" "  -0.63%  "
Set-TextValue "D32" "0.999"
Set-TextValue "E32" "  -0.11%  "
Set-TextValue "D33" "18.18"
Set-TextValue "E33" "  +6.34%  "
Set-TextValue "E34" "  -0.85%  "
Set-TextValue "D35" "4.50"
Set-TextValue "E35" "  -8.06%  "
Set-TextValue "D36" "2.34"
Set-TextValue "E36" "  -1.39%  "
Set-TextValue "E37" "  +0.52%  "
Set-TextValue "E38" "  -0.57%  "
Set-TextValue "E39" "  -0.42%  "
Set-TextValue "E40" "  +1.32%  "
Set-TextValue "E41" "  -0.72%  "
Set-TextValue "D42" "1.996.43"
Set-TextValue "E43" "  -0.49%  "
Set-TextValue "D44" "10.29"
Set-TextValue "E44" "  +1.16%  "
Set-TextValue "B45" "EnergySwap"
Set-TextValue "C45" "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue "D45" "18.09"
Set-TextValue "E45" "  +4.04%  "
Set-TextValue "B46" "ApeXProtocol"
Set-TextValue "C46" "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
Set-TextValue "D46" "2.11"
Set-TextValue "E46" "  -1.40%  "
Set-TextValue "D47" "2.77"
Set-TextValue "E47" "  -0.85%  "
Set-TextValue "D48" "2.531.09"
Set-TextValue "E48" "  +0.36%  "
Set-TextValue "B49" "HuobiToken"
Set-TextValue "C49" "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
Set-TextValue "D49" "2.84"
Set-TextValue "E49" "  -3.40%  "
Set-TextValue "B50" "MultiversX"
Set-TextValue "C50" "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
Set-TextValue "D50" "53.40"
Set-TextValue "E50" "  -0.20%  "
Set-TextValue "D51" "71.76"
Set-TextValue "E51" "  -0.45%  "

Write-Output "Applied 89 cell updates"
